# Weekly update: prepend two new price records (week of 2023-12-07)
# for "Agrícola del Norte S.A. de Arica - Plátano" and push the
# existing history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (row 440),
# shifting the existing rows (440..456) down to (442..458).
$ws.Rows("440:441").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# New row 440: "Pintón" record
$ws.Cells.Item(440, 1).Value2  = 1
$ws.Cells.Item(440, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(440, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(440, 4).Value2  = 45267
$ws.Cells.Item(440, 5).Value2  = 15
$ws.Cells.Item(440, 6).Value2  = "Fruta"
$ws.Cells.Item(440, 7).Value2  = 100108
$ws.Cells.Item(440, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(440, 9).Value2  = 100108006
$ws.Cells.Item(440, 10).Value2 = "Plátano"
$ws.Cells.Item(440, 11).Value2 = "Sin especificar"
$ws.Cells.Item(440, 12).Value2 = "Pintón"
$ws.Cells.Item(440, 13).Value2 = 108
$ws.Cells.Item(440, 14).Value2 = 20000
$ws.Cells.Item(440, 15).Value2 = 21000
$ws.Cells.Item(440, 16).Value2 = 20500
$ws.Cells.Item(440, 17).Value2 = "`$/caja 20 kilos"
$ws.Cells.Item(440, 18).Value2 = "Ecuador"
$ws.Cells.Item(440, 19).Value2 = 1025
$ws.Cells.Item(440, 20).Value2 = 20

# New row 441: "Verde" record
$ws.Cells.Item(441, 1).Value2  = 1
$ws.Cells.Item(441, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(441, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(441, 4).Value2  = 45267
$ws.Cells.Item(441, 5).Value2  = 15
$ws.Cells.Item(441, 6).Value2  = "Fruta"
$ws.Cells.Item(441, 7).Value2  = 100108
$ws.Cells.Item(441, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(441, 9).Value2  = 100108006
$ws.Cells.Item(441, 10).Value2 = "Plátano"
$ws.Cells.Item(441, 11).Value2 = "Sin especificar"
$ws.Cells.Item(441, 12).Value2 = "Verde"
$ws.Cells.Item(441, 13).Value2 = 108
$ws.Cells.Item(441, 14).Value2 = 21000
$ws.Cells.Item(441, 15).Value2 = 22000
$ws.Cells.Item(441, 16).Value2 = 21500
$ws.Cells.Item(441, 17).Value2 = "`$/caja 20 kilos"
$ws.Cells.Item(441, 18).Value2 = "Ecuador"
$ws.Cells.Item(441, 19).Value2 = 1075
$ws.Cells.Item(441, 20).Value2 = 20
